$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "test1"
$ws.Range("A3").Value = "test32"

$ws.Range("B7").Select()
